# Fruta / hortaliza, semanal
# A new weekly price observation (Ajo, Chino, China) is inserted as row 123,
# pushing all existing rows from 123 onward down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 123; this shifts the former rows
# 123..193 down to 124..194 (and grows the sheet's dimension accordingly).
$ws.Rows("123:123").Insert()

# Populate the newly inserted row with the new observation's data.
$ws.Range("A123").Value = 5
$ws.Range("B123").Value = "Macroferia Regional de Talca"
$ws.Range("C123").Value = "Maule"
$ws.Range("D123").Value = 44518
$ws.Range("E123").Value = 7
$ws.Range("F123").Value = 100112003
$ws.Range("G123").Value = "Ajo"
$ws.Range("H123").Value = "Chino"
$ws.Range("I123").Value = "Primera"
$ws.Range("J123").Value = 200
$ws.Range("K123").Value = 20000
$ws.Range("L123").Value = 20000
$ws.Range("M123").Value = 20000
$ws.Range("N123").Value = "`$/caja 10 kilos"
$ws.Range("O123").Value = "China"
$ws.Range("P123").Value = 2000
$ws.Range("Q123").Value = 10
$ws.Range("R123").Value = "Hortaliza"
